# Scheduled runner update: refresh Leve-profit market-board figures across
# all Balmung_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Generated from the upstream commit's per-cell value changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2522.7144
$ws.Range("J40").Value = 2132.8
$ws.Range("L40").Value = 2132.8
$ws.Range("N40").Value = -2482.8
$ws.Range("H45").Value = 18140.75
$ws.Range("I45").Value = 1047.6666
$ws.Range("K45").Value = 3142.9998
$ws.Range("M45").Value = -2950.9998
$ws.Range("H58").Value = 204.33333
$ws.Range("I58").Value = 358.66666
$ws.Range("J58").Value = 50
$ws.Range("K58").Value = 1075.99998
$ws.Range("L58").Value = 150
$ws.Range("M58").Value = -925.9999800000001
$ws.Range("N58").Value = -450
$ws.Range("H70").Value = 831.7143
$ws.Range("J70").Value = 861.3333
$ws.Range("L70").Value = 2583.9999
$ws.Range("N70").Value = -3123.9999
$ws.Range("H73").Value = 831.7143
$ws.Range("J73").Value = 861.3333
$ws.Range("L73").Value = 2583.9999
$ws.Range("N73").Value = -4455.9999
$ws.Range("H74").Value = 12499.75
$ws.Range("I74").Value = 12499.75
$ws.Range("K74").Value = 12499.75
$ws.Range("M74").Value = -11563.75
$ws.Range("H77").Value = 12499.75
$ws.Range("I77").Value = 12499.75
$ws.Range("K77").Value = 62498.75
$ws.Range("M77").Value = -57818.75
$ws.Range("H107").Value = 1171.2354
$ws.Range("I107").Value = 1348.4286
$ws.Range("K107").Value = 1348.4286
$ws.Range("M107").Value = 571.5714
$ws.Range("H112").Value = 77087.64
$ws.Range("J112").Value = 82710.234
$ws.Range("L112").Value = 248130.702
$ws.Range("N112").Value = -250346.702
$ws.Range("H138").Value = 12766
$ws.Range("I138").Value = 16012.4
$ws.Range("J138").Value = 4650
$ws.Range("K138").Value = 48037.2
$ws.Range("L138").Value = 13950
$ws.Range("M138").Value = -42897.2
$ws.Range("N138").Value = -24230
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 149514.03
$ws.Range("I32").Value = 176142.45
$ws.Range("J32").Value = 11530.363
$ws.Range("K32").Value = 176142.45
$ws.Range("L32").Value = 11530.363
$ws.Range("M32").Value = -175855.45
$ws.Range("N32").Value = -12104.363
$ws.Range("H61").Value = 1715953.6
$ws.Range("I61").Value = 69831.19
$ws.Range("K61").Value = 69831.19
$ws.Range("M61").Value = -69619.19
$ws.Range("H132").Value = 3110.2354
$ws.Range("I132").Value = 3279
$ws.Range("K132").Value = 9837
$ws.Range("M132").Value = -7307
$ws.Range("H136").Value = 1715953.6
$ws.Range("I136").Value = 69831.19
$ws.Range("K136").Value = 209493.57
$ws.Range("M136").Value = -206943.57
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2197.8572
$ws.Range("I5").Value = 914.1667
$ws.Range("J5").Value = 9900
$ws.Range("K5").Value = 914.1667
$ws.Range("L5").Value = 9900
$ws.Range("M5").Value = -801.1667
$ws.Range("N5").Value = -10126
$ws.Range("H86").Value = 3560
$ws.Range("I86").Value = 2734.4583
$ws.Range("K86").Value = 2734.4583
$ws.Range("M86").Value = -1611.4583
$ws.Range("H89").Value = 3560
$ws.Range("I89").Value = 2734.4583
$ws.Range("K89").Value = 13672.2915
$ws.Range("M89").Value = -8056.291499999999
$ws.Range("H134").Value = 25716168
$ws.Range("I134").Value = 2041.1666
$ws.Range("K134").Value = 6123.4998
$ws.Range("M134").Value = -3588.4998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 583.5217
$ws.Range("J22").Value = 877
$ws.Range("L22").Value = 877
$ws.Range("N22").Value = -1577
$ws.Range("H31").Value = 6867.926
$ws.Range("I31").Value = 5302.3887
$ws.Range("K31").Value = 5302.3887
$ws.Range("M31").Value = -5007.3887
$ws.Range("H34").Value = 6867.926
$ws.Range("I34").Value = 5302.3887
$ws.Range("K34").Value = 5302.3887
$ws.Range("M34").Value = -5100.3887
$ws.Range("H132").Value = 2600
$ws.Range("I132").Value = 1772.8462
$ws.Range("K132").Value = 5318.5386
$ws.Range("M132").Value = -2788.5386
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9133
$ws.Range("J32").Value = 9949.5
$ws.Range("L32").Value = 29848.5
$ws.Range("N32").Value = -30414.5
$ws.Range("H46").Value = 1186.875
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3182
$ws.Range("H119").Value = 100013000
$ws.Range("I119").Value = 166672340
$ws.Range("K119").Value = 500017020
$ws.Range("M119").Value = -500012182
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 817599.8
$ws.Range("J11").Value = 1728001.1
$ws.Range("L11").Value = 1728001.1
$ws.Range("N11").Value = -1728279.1
$ws.Range("H80").Value = 7898834.5
$ws.Range("I80").Value = 95686.336
$ws.Range("K80").Value = 95686.336
$ws.Range("M80").Value = -94688.336
$ws.Range("H83").Value = 7898834.5
$ws.Range("I83").Value = 95686.336
$ws.Range("K83").Value = 478431.68
$ws.Range("M83").Value = -473439.68
$ws.Range("H104").Value = 54835.5
$ws.Range("J104").Value = 54835.5
$ws.Range("L104").Value = 54835.5
$ws.Range("N104").Value = -61823.5
$ws.Range("H132").Value = 12821737
$ws.Range("I132").Value = 2450
$ws.Range("K132").Value = 7350
$ws.Range("M132").Value = -4820
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3905.7222
$ws.Range("I22").Value = 695.1429000000001
$ws.Range("J22").Value = 5948.8184
$ws.Range("K22").Value = 695.1429000000001
$ws.Range("L22").Value = 5948.8184
$ws.Range("M22").Value = -400.1429000000001
$ws.Range("N22").Value = -6538.8184
$ws.Range("H27").Value = 3905.7222
$ws.Range("I27").Value = 695.1429000000001
$ws.Range("J27").Value = 5948.8184
$ws.Range("K27").Value = 695.1429000000001
$ws.Range("L27").Value = 5948.8184
$ws.Range("M27").Value = -588.1429000000001
$ws.Range("N27").Value = -6162.8184
$ws.Range("H46").Value = 8075.8945
$ws.Range("I46").Value = 35659
$ws.Range("J46").Value = 2904.0625
$ws.Range("K46").Value = 35659
$ws.Range("L46").Value = 2904.0625
$ws.Range("M46").Value = -35471
$ws.Range("N46").Value = -3280.0625
$ws.Range("H68").Value = 3334.5
$ws.Range("I68").Value = 2694.1538
$ws.Range("K68").Value = 2694.1538
$ws.Range("M68").Value = -1945.1538
$ws.Range("H71").Value = 3334.5
$ws.Range("I71").Value = 2694.1538
$ws.Range("K71").Value = 13470.769
$ws.Range("M71").Value = -9726.769
$ws.Range("H122").Value = 3419.5
$ws.Range("J122").Value = 3921.25
$ws.Range("L122").Value = 11763.75
$ws.Range("N122").Value = -16663.75
$ws.Range("H132").Value = 6990
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 7250
$ws.Range("I52").Value = 7250
$ws.Range("K52").Value = 7250
$ws.Range("M52").Value = -7024
$ws.Range("H81").Value = 75756.28999999999
$ws.Range("I81").Value = 4880.8887
$ws.Range("J81").Value = 203332
$ws.Range("K81").Value = 9761.777400000001
$ws.Range("L81").Value = 406664
$ws.Range("M81").Value = -8700.777400000001
$ws.Range("N81").Value = -408786
$ws.Range("H84").Value = 75756.28999999999
$ws.Range("I84").Value = 4880.8887
$ws.Range("J84").Value = 203332
$ws.Range("K84").Value = 48808.887
$ws.Range("L84").Value = 2033320
$ws.Range("M84").Value = -43504.887
$ws.Range("N84").Value = -2043928
$ws.Range("H132").Value = 1918.2632
$ws.Range("I132").Value = 1722.2903
$ws.Range("J132").Value = 2786.1428
$ws.Range("K132").Value = 5166.8709
$ws.Range("L132").Value = 8358.428400000001
$ws.Range("M132").Value = -2636.8709
$ws.Range("N132").Value = -13418.4284
